$wb = $excel.ActiveWorkbook

# --- Sheet1 -> TestData -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TestData"

# Add a new "Address" column (G) to the third test block. Column G already
# carries the sheet's default (centered) column style, so the new cells
# automatically pick up the same style used by the rest of the row.
$ws1.Range("G12").Value = "Address"
$ws1.Range("G13").Value = "Thane"
$ws1.Range("G14").Value = "Thane"
$ws1.Range("G15").Value = "Thane"
$ws1.Range("G16").Value = "Thane"
$ws1.Range("G17").Value = "Thane"
$ws1.Range("G18").Value = "Thane"

# Rename the three test-case headers (Radio/Telephony/Hybrid -> TestA/TestB/TestC)
$ws1.Range("A1").Value = "TestA"
$ws1.Range("A5").Value = "TestB"
$ws1.Range("A11").Value = "TestC"

# --- New sheet: TestCaseList --------------------------------------------
$notes = $wb.Worksheets.Item("Notes")
$ws3 = $wb.Worksheets.Add($null, $notes)
$ws3.Name = "TestCaseList"

$ws3.Range("B1").Value = "ExecutionMode"
$ws3.Range("A1").Value = "TesCaseName"
$ws3.Range("A2").Value = "TestA"
$ws3.Range("B2").Value = "Yes"
$ws3.Range("A3").Value = "TestB"
$ws3.Range("B3").Value = "No"
$ws3.Range("A4").Value = "TestC"
$ws3.Range("B4").Value = "No"

# TestCaseList becomes the active/visible tab
$ws3.Activate()
